$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2.55802675585284
$ws.Range("C2").Value = 2.82095238095238
$ws.Range("D2").Value = 2.6030888030888
$ws.Range("E2").Value = 3.00679611650485

$ws.Range("B3").Value = 4.37597402597403
$ws.Range("C3").Value = 4.55853658536585
$ws.Range("D3").Value = 3.15253623188406
$ws.Range("E3").Value = 2.67467532467532

$ws.Range("B4").Value = 3.12065972222222
$ws.Range("C4").Value = 2.99878048780488
$ws.Range("D4").Value = 1.89588235294118
$ws.Range("E4").Value = 2.01393939393939

$ws.Range("B5").Value = 3.76937919463087
$ws.Range("C5").Value = 3.77865853658537
$ws.Range("D5").Value = 2.54905838041431
$ws.Range("E5").Value = 2.33291536050157
